$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "How many different lithology types can be recorded in a log at most?"
$ws.Range("B46").Value = "The maximum number of lithology types that can be recorded in a log is 450."
